$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.036642952388835
$ws.Range("D2").Value = 1.039738277539773
$ws.Range("E2").Value = 1.040291862911117
$ws.Range("F2").Value = 1.035332380744921
$ws.Range("I2").Value = 1.037613750980864
$ws.Range("J2").Value = 1.04175015313205
$ws.Range("K2").Value = 1.042522612503839
$ws.Range("L2").Value = 1.043074626300238
$ws.Range("M2").Value = 1.038129288691706
$ws.Range("N2").Value = 1.017745021418722
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.037863676982625
$ws.Range("D3").Value = 1.040681248606848
$ws.Range("E3").Value = 1.041461484360525
$ws.Range("F3").Value = 1.037184107460091
$ws.Range("I3").Value = 1.037967134636404
$ws.Range("J3").Value = 1.042613562020152
$ws.Range("K3").Value = 1.043275514051108
$ws.Range("L3").Value = 1.044053696011068
$ws.Range("M3").Value = 1.039787619361562
$ws.Range("N3").Value = 1.018041305644086
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038652500384664
$ws.Range("D4").Value = 1.041290447999952
$ws.Range("E4").Value = 1.04221760150363
$ws.Range("F4").Value = 1.038381021438111
$ws.Range("I4").Value = 1.038194078464145
$ws.Range("J4").Value = 1.043170687962958
$ws.Range("K4").Value = 1.043761106910434
$ws.Range("L4").Value = 1.044685943549353
$ws.Range("M4").Value = 1.040858979652045
$ws.Range("N4").Value = 1.018232261827174
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.038983870691209
$ws.Range("D5").Value = 1.041546326850872
$ws.Range("E5").Value = 1.042535307802231
$ws.Range("F5").Value = 1.038883910036971
$ws.Range("I5").Value = 1.038289075575947
$ws.Range("J5").Value = 1.043404534141767
$ws.Range("K5").Value = 1.043964873409088
$ws.Range("L5").Value = 1.044951438370533
$ws.Range("M5").Value = 1.041308985224976
$ws.Range("N5").Value = 1.018312358934858
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039039494629369
$ws.Range("D6").Value = 1.041589276709495
$ws.Range("E6").Value = 1.042588642498831
$ws.Range("F6").Value = 1.038968330341428
$ws.Range("I6").Value = 1.038305002004506
$ws.Range("J6").Value = 1.043443776366964
$ws.Range("K6").Value = 1.043999064685632
$ws.Range("L6").Value = 1.04499599850681
$ws.Range("M6").Value = 1.04138452032115
$ws.Range("N6").Value = 1.018325797013009
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038656929151932
$ws.Range("D7").Value = 1.041293867960523
$ws.Range("E7").Value = 1.042221847355866
$ws.Range("F7").Value = 1.038387742201108
$ws.Range("I7").Value = 1.038195349429705
$ws.Range("J7").Value = 1.043173814075165
$ws.Range("K7").Value = 1.043763831125409
$ws.Range("L7").Value = 1.044689492286724
$ws.Range("M7").Value = 1.040864994184964
$ws.Range("N7").Value = 1.018233332798166
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037055725573397
$ws.Range("D8").Value = 1.040057160647984
$ws.Range("E8").Value = 1.040687290404992
$ws.Range("F8").Value = 1.035958451926596
$ws.Range("I8").Value = 1.037733535463006
$ws.Range("J8").Value = 1.042042270637085
$ws.Range("K8").Value = 1.042777388855474
$ws.Range("L8").Value = 1.043405773782168
$ws.Range("M8").Value = 1.038690085635393
$ws.Range("N8").Value = 1.017845309904221
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.034225830044831
$ws.Range("D9").Value = 1.037870405515538
$ws.Range("E9").Value = 1.037977629801778
$ws.Range("F9").Value = 1.031667426823777
$ws.Range("I9").Value = 1.036906528392526
$ws.Range("J9").Value = 1.040036272047263
$ws.Range("K9").Value = 1.041026893837978
$ws.Range("L9").Value = 1.041133769764197
$ws.Range("M9").Value = 1.034844203091048
$ws.Range("N9").Value = 1.01715570354114
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.032333325969597
$ws.Range("D10").Value = 1.036407341673452
$ws.Range("E10").Value = 1.036167204746022
$ws.Range("F10").Value = 1.028799079742597
$ws.Range("I10").Value = 1.036346200070086
$ws.Range("J10").Value = 1.038690613471609
$ws.Range("K10").Value = 1.039851485495099
$ws.Range("L10").Value = 1.039612206023919
$ws.Range("M10").Value = 1.032270605199295
$ws.Range("N10").Value = 1.016691963932852
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.031512386805493
$ws.Range("D11").Value = 1.035772539978639
$ws.Range("E11").Value = 1.03538227220481
$ws.Range("F11").Value = 1.027555072168239
$ws.Range("I11").Value = 1.036101418001083
$ws.Range("J11").Value = 1.038105906604215
$ws.Range("K11").Value = 1.039340486971011
$ws.Range("L11").Value = 1.03895166729156
$ws.Range("M11").Value = 1.031153767426434
$ws.Range("N11").Value = 1.016490196242529
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.031207226717638
$ws.Range("D12").Value = 1.035536549779873
$ws.Range("E12").Value = 1.035090557390904
$ws.Range("F12").Value = 1.027092678567816
$ws.Range("I12").Value = 1.036010169191308
$ws.Range("J12").Value = 1.037888411642597
$ws.Range("K12").Value = 1.03915036962049
$ws.Range("L12").Value = 1.038706055288028
$ws.Range("M12").Value = 1.030738542958027
$ws.Range("N12").Value = 1.016415104444388
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.031272694941427
$ws.Range("D13").Value = 1.035587179430528
$ws.Range("E13").Value = 1.035153138340249
$ws.Range("F13").Value = 1.027191877993428
$ws.Range("I13").Value = 1.036029757146933
$ws.Range("J13").Value = 1.037935079090725
$ws.Range("K13").Value = 1.039191164509413
$ws.Range("L13").Value = 1.038758751665052
$ws.Range("M13").Value = 1.030827627470288
$ws.Range("N13").Value = 1.016431218530435
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.031487166834727
$ws.Range("D14").Value = 1.035753036985649
$ws.Range("E14").Value = 1.035358162175307
$ws.Range("F14").Value = 1.027516857052728
$ws.Range("I14").Value = 1.036093882004329
$ws.Range("J14").Value = 1.038087934726278
$ws.Range("K14").Value = 1.03932477815012
$ws.Range("L14").Value = 1.038931370225598
$ws.Range("M14").Value = 1.031119452692217
$ws.Range("N14").Value = 1.016483992122744
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.031619279836951
$ws.Range("D15").Value = 1.035855201145608
$ws.Range("E15").Value = 1.03548446332474
$ws.Range("F15").Value = 1.027717045369778
$ws.Range("I15").Value = 1.03613334820975
$ws.Range("J15").Value = 1.038182073085172
$ws.Range("K15").Value = 1.03940706077681
$ws.Range("L15").Value = 1.039037691831846
$ws.Range("M15").Value = 1.031299204950399
$ws.Range("N15").Value = 1.016516488251486
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.032387777506305
$ws.Range("D16").Value = 1.036449443994174
$ws.Range("E16").Value = 1.03621927652694
$ws.Range("F16").Value = 1.028881597202671
$ws.Range("I16").Value = 1.036362399858833
$ws.Range("J16").Value = 1.038729375464289
$ws.Range("K16").Value = 1.039885355534863
$ws.Range("L16").Value = 1.039656007792517
$ws.Range("M16").Value = 1.032344673284942
$ws.Range("N16").Value = 1.016705334143528
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.032869437702934
$ws.Range("D17").Value = 1.036821850257979
$ws.Range("E17").Value = 1.036679932650845
$ws.Range("F17").Value = 1.029611545886568
$ws.Range("I17").Value = 1.036505499289274
$ws.Range("J17").Value = 1.039072137969659
$ws.Range("K17").Value = 1.040184828911507
$ws.Range("L17").Value = 1.04004340510142
$ws.Range("M17").Value = 1.032999802324124
$ws.Range("N17").Value = 1.016823532839405
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.033150240255976
$ws.Range("D18").Value = 1.037038944745858
$ws.Range("E18").Value = 1.036948528777478
$ws.Range("F18").Value = 1.030037121140132
$ws.Range("I18").Value = 1.036588758853144
$ws.Range("J18").Value = 1.039271870166465
$ws.Range("K18").Value = 1.040359310168884
$ws.Range("L18").Value = 1.040269204480932
$ws.Range("M18").Value = 1.033381692200452
$ws.Range("N18").Value = 1.016892383079686
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.033245962802696
$ws.Range("D19").Value = 1.037112947463953
$ws.Range("E19").Value = 1.037040096838226
$ws.Range("F19").Value = 1.03018219925494
$ws.Range("I19").Value = 1.036617113004374
$ws.Range("J19").Value = 1.039339940669703
$ws.Range("K19").Value = 1.04041877054204
$ws.Range("L19").Value = 1.040346168731874
$ws.Range("M19").Value = 1.033511867163877
$ws.Range("N19").Value = 1.016915843485609
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.032817774819299
$ws.Range("D20").Value = 1.036781907420133
$ws.Range("E20").Value = 1.036630518655789
$ws.Range("F20").Value = 1.02953324921005
$ws.Range("I20").Value = 1.036490167594411
$ws.Range("J20").Value = 1.039035383036158
$ws.Range("K20").Value = 1.040152718602629
$ws.Range("L20").Value = 1.040001857920458
$ws.Range("M20").Value = 1.032929537663103
$ws.Range("N20").Value = 1.016810860872991
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.031424016511099
$ws.Range("D21").Value = 1.035704201524284
$ws.Range("E21").Value = 1.035297792108018
$ws.Range("F21").Value = 1.027421167596354
$ws.Range("I21").Value = 1.036075007842542
$ws.Range("J21").Value = 1.038042931109162
$ws.Range("K21").Value = 1.039285440841275
$ws.Range("L21").Value = 1.038880545533248
$ws.Range("M21").Value = 1.03103352806295
$ws.Range("N21").Value = 1.016468455663752
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.030546390791541
$ws.Range("D22").Value = 1.035025465932148
$ws.Range("E22").Value = 1.034458949577079
$ws.Range("F22").Value = 1.026091395421553
$ws.Range("I22").Value = 1.035812094203495
$ws.Range("J22").Value = 1.03741714851243
$ws.Range("K22").Value = 1.038738355078544
$ws.Range("L22").Value = 1.038174034094882
$ws.Range("M22").Value = 1.029839219765174
$ws.Range("N22").Value = 1.016252324873489
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.031011763247647
$ws.Range("D23").Value = 1.035385385653617
$ws.Range("E23").Value = 1.034903723184135
$ws.Range("F23").Value = 1.026796510335344
$ws.Range("I23").Value = 1.035951649130397
$ws.Range("J23").Value = 1.037749058644144
$ws.Range("K23").Value = 1.039028546773545
$ws.Range("L23").Value = 1.038548712620227
$ws.Range("M23").Value = 1.030472559101098
$ws.Range("N23").Value = 1.016366980607635
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.032841119490046
$ws.Range("D24").Value = 1.036799956251528
$ws.Range("E24").Value = 1.036652847011957
$ws.Range("F24").Value = 1.02956862869817
$ws.Range("I24").Value = 1.036497095969671
$ws.Range("J24").Value = 1.039051991611901
$ws.Range("K24").Value = 1.040167228475965
$ws.Range("L24").Value = 1.040020631806259
$ws.Range("M24").Value = 1.032961287966321
$ws.Range("N24").Value = 1.016816587076623
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.034958446422227
$ws.Range("D25").Value = 1.038436642739641
$ws.Range("E25").Value = 1.038678826902359
$ws.Range("F25").Value = 1.032778056244497
$ws.Range("I25").Value = 1.037121907204055
$ws.Range("J25").Value = 1.040556322957963
$ws.Range("K25").Value = 1.041480908794401
$ws.Range("L25").Value = 1.041722336251369
$ws.Range("M25").Value = 1.03584011302242
$ws.Range("N25").Value = 1.017334684130749
